# Implemented getting standard relationship between microservices and started
# implementing MSM measure.
#
# This updates the "classFields" sheet so that several field records line up
# with the correct class/field groupings (e.g. the KeycloakRole enum's
# BLOGGER/$VALUES rows were swapped, and the allowed*/exposed* CORS field
# rows for CrossOriginRequestSharingFilter / SecurityConfiguration /
# KeycloakRealmRoleConverter were in the wrong order).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Row 2 <-> Row 4 (org.andante.config.security.role.KeycloakRole fields):
# swap BLOGGER/public/KeycloakRole with $VALUES/private/KeycloakRole[]
$ws.Range("B2").Value = "`$VALUES"
$ws.Range("C2").Value = "private"
$ws.Range("D2").Value = "org.andante.config.security.role.KeycloakRole[]"

$ws.Range("B4").Value = "BLOGGER"
$ws.Range("C4").Value = "public"
$ws.Range("D4").Value = "org.andante.config.security.role.KeycloakRole"

# Rows 7-10 (org.andante.config.security.filter.CrossOriginRequestSharingFilter):
# swap allowedHeaders <-> allowedOrigins, allowedMethods <-> exposedHeaders
$ws.Range("B7").Value = "allowedOrigins"
$ws.Range("B8").Value = "exposedHeaders"
$ws.Range("B9").Value = "allowedHeaders"
$ws.Range("B10").Value = "allowedMethods"

# Rows 11-12 (org.andante.config.security.converter.KeycloakRealmRoleConverter):
# swap REALM_ACCESS <-> ROLES
$ws.Range("B11").Value = "ROLES"
$ws.Range("B12").Value = "REALM_ACCESS"

# Row 16 <-> Row 19 (org.andante.config.security.SecurityConfiguration):
# swap jwkSetUri/java.lang.String with allowedMethods/java.util.List
$ws.Range("B16").Value = "allowedMethods"
$ws.Range("D16").Value = "java.util.List"

$ws.Range("B19").Value = "jwkSetUri"
$ws.Range("D19").Value = "java.lang.String"
